$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.441.80'
$ws.Range("E2").Value = '  -2.14%  '
$ws.Range("D3").Value = '3.025.95'
$ws.Range("E3").Value = '  -4.63%  '
$ws.Range("E4").Value = '  -0.18%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '567.67'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -3.79%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '130.03'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -4.96%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.026.54'
$ws.Range("E8").Value = '  -4.38%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.499'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -2.43%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.136'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -4.89%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '5.28'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +0.12%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.435'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -4.96%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.0000225'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -4.44%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '33.26'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -4.99%  '
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").Value = '3.509.39'
$ws.Range("E16").Value = '  -5.00%  '
$ws.Range("D17").Value = '61.351.25'
$ws.Range("E17").Value = '  -2.27%  '
$ws.Range("D18").Value = '3.009.24'
$ws.Range("E18").Value = '  -5.14%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '6.26'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -5.12%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '440.46'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -3.85%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '13.26'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -4.95%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.668'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -5.99%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '7.19'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -5.78%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '12.80'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -4.54%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '79.40'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -4.80%  '
$ws.Range("E26").Value = '  -0.06%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.14%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '2.51'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -6.31%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '7.30'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -6.03%  '
$ws.Range("E30").Value = '  -6.37%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '6.22'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -9.07%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '25.70'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -6.09%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.0947'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -8.54%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '2.30'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -3.39%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.961'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -7.53%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '5.61'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -3.91%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '50.46'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("D38").Value = '0.0₃0679'
$ws.Range("E38").Value = '  -3.76%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.0364'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -6.23%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '7.82'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -3.56%  '
$ws.Range("E41").Value = '  -2.17%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '380.36'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -5.19%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '2.49'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -7.52%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.684.37'
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("B46").Value = 'Arweave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '34.86'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.238'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -5.88%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '2.00'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -6.18%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '119.34'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -4.75%  '
$ws.Range("E50").Value = '  -3.62%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '23.61'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -7.74%  '
